$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $text) {
    # Force the value to be written as text (matching the original inlineStr
    # cells), even when the string looks like a number, without leaving a
    # lingering custom number-format style on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Column D (Price) and Column E (Volume(1h)) updates, keyed by row number.
$updates = @{
    2  = @{ D = "20.101.04"; E = "  -1.74%  " }
    3  = @{ D = "1.425.63";  E = "  -1.46%  " }
    4  = @{ D = "0.9972";    E = "  -0.77%  " }
    5  = @{ D = "0.9957";    E = "  -0.96%  " }
    6  = @{ D = "276.71";    E = "  -0.49%  " }
    7  = @{ D = "0.3701";    E = "  -0.53%  " }
    8  = @{ D = "0.3132";    E = "  +2.07%  " }
    9  = @{ D = "39.68";     E = "  -2.84%  " }
    10 = @{ D = "1.056";     E = "  +4.36%  " }
    11 = @{ D = "0.06551";   E = "  +0.08%  " }
    12 = @{ D = "0.9956";    E = "  -0.92%  " }
    13 = @{ D = "5.520";     E = "  +2.41%  " }
    14 = @{ D = "17.83";     E = "  +3.65%  " }
    15 = @{ D = "6.215";     E = "  +1.23%  " }
    16 = @{ D = "1.421.92";  E = "  -1.78%  " }
    17 = @{ D = "0.00001024"; E = "  +0.86%  " }
    18 = @{ D = "0.05686";   E = "  -3.21%  " }
    19 = @{ D = "0.9960";    E = "  -1.02%  " }
    20 = @{ D = "71.52";     E = "  -6.60%  " }
    21 = @{ D = "5.623";     E = "  -1.73%  " }
    22 = @{ D = "14.89";     E = "  +3.47%  " }
    23 = @{ D = "11.08";     E = "  +1.68%  " }
    24 = @{ D = "2.243";     E = "  -1.66%  " }
    25 = @{ D = "20.101.28"; E = "  -1.65%  " }
    26 = @{ D = "2.295";     E = "  +3.33%  " }
    27 = @{ D = "133.82";    E = "  -6.58%  " }
    28 = @{ D = "17.30";     E = "  +1.63%  " }
    29 = @{ D = "1.578.64";  E = "  -2.06%  " }
    30 = @{ D = "110.73";    E = "  +1.32%  " }
    31 = @{ D = "3.923";     E = "  +6.76%  " }
    32 = @{ D = "5.284";     E = "  -2.16%  " }
    33 = @{ D = "0.8266";    E = "  -9.38%  " }
    34 = @{ D = "0.07805";   E = "  +0.91%  " }
    35 = @{ D = "1.466";     E = "  +4.11%  " }
    36 = @{ D = "4.937";     E = "  +4.28%  " }
    37 = @{ D = "0.05869";   E = "  +4.75%  " }
    38 = @{ E = "  -2.77%  " }
    39 = @{ D = "0.9955";    E = "  -1.00%  " }
    42 = @{ D = "1.109";     E = "  -1.36%  " }
    43 = @{ D = "0.1882";    E = "  -1.63%  " }
    44 = @{ D = "0.5343";    E = "  +0.42%  " }
    45 = @{ D = "12.38";     E = "  +2.60%  " }
    46 = @{ D = "3.542";     E = "  -1.21%  " }
    47 = @{ D = "117.78";    E = "  +5.86%  " }
    48 = @{ D = "0.5232";    E = "  +1.50%  " }
    49 = @{ D = "1.783";     E = "  +0.31%  " }
    50 = @{ D = "1.040";     E = "  -1.51%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $vals.D
    }
    if ($vals.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$row") $vals.E
    }
}

# Rows 40/41 swap: VeChain <-> Aptos (Aptos now ranks ahead of VeChain)
Set-TextValue $ws.Range("B40") "Aptos"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "10.64"
Set-TextValue $ws.Range("E40") "  -1.21%  "

Set-TextValue $ws.Range("B41") "VeChain"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.02064"
Set-TextValue $ws.Range("E41") "  +1.24%  "

# Row 51: PaxDollar -> Cronos
Set-TextValue $ws.Range("B51") "Cronos"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.06234"
Set-TextValue $ws.Range("E51") "  -1.00%  "
